# The deck currently uses the "Integral" / "Red Violet" design (Design tab)
# for every slide. This edit applies the built-in "Office Theme" design in
# its place - i.e. the same operation as opening the Design tab in
# PowerPoint and clicking the plain "Office Theme" thumbnail.
#
# That swaps out the theme colour palette (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink) that is wired to the slide master's theme; the font
# scheme and format scheme are identical between the two built-in themes,
# so only the twelve theme colours actually change.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# All slides share one Design/SlideMaster, so any slide's theme colour
# scheme is the live, shared one for the whole deck.
$s = $p.Slides.Item(1)
$scheme = $s.ThemeColorScheme

# Target palette: the stock PowerPoint "Office Theme" colours, replacing
# the current "Red Violet" ("Integral" design) colours.
$officeColors = @(
    (RGB 0x00 0x00 0x00),  # 1  dk1
    (RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (RGB 0x44 0x54 0x6A),  # 3  dk2
    (RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (RGB 0xED 0x7D 0x31),  # 6  accent2
    (RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (RGB 0xFF 0xC0 0x00),  # 8  accent4
    (RGB 0x44 0x72 0xC4),  # 9  accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}

Write-Host "Applied Office Theme colours to the active design."
